$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.464.47"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "1.878.15"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'247.01"
$ws.Range("E5").Value = "  +5.58%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.4760"
$ws.Range("E7").Value = "  +2.12%  "
$ws.Range("D8").Value = "'0.2901"
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("D9").Value = "'0.06522"
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").Value = "'21.58"
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("D11").Value = "'0.07735"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "'0.7452"
$ws.Range("E12").Value = "  +9.25%  "
$ws.Range("D13").Value = "'97.00"
$ws.Range("E13").Value = "  +3.59%  "
$ws.Range("D14").Value = "1.879.14"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "'5.119"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "'273.93"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "30.447.35"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "'13.62"
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("D19").Value = "'0.000007569"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "2.125.87"
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'5.248"
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("D24").Value = "'6.170"
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("D25").Value = "'9.287"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("E27").Value = "  +1.78%  "
$ws.Range("D28").Value = "'1.958"
$ws.Range("E28").Value = "  +3.68%  "
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").Value = "'0.09996"
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("D31").Value = "'1.516"
$ws.Range("E31").Value = "  +4.53%  "
$ws.Range("D32").Value = "'4.326"
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("D34").Value = "'0.04774"
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("D35").Value = "'1.124"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").Value = "'0.6992"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D37").Value = "'2.716"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").Value = "'0.01867"
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").Value = "'6.358"
$ws.Range("E40").Value = "  +0.83%  "
$ws.Range("D41").Value = "'1.932"
$ws.Range("E41").Value = "  +2.80%  "
$ws.Range("D42").Value = "'70.09"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("D43").Value = "'0.4169"
$ws.Range("E43").Value = "  +2.78%  "
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'0.8385"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").Value = "'102.68"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("D47").Value = "'9.343"
$ws.Range("E47").Value = "  +3.84%  "
$ws.Range("E48").Value = "  +1.92%  "
$ws.Range("D49").Value = "'35.33"
$ws.Range("E49").Value = "  +3.90%  "
$ws.Range("D50").Value = "'926.21"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").Value = "'0.05610"
